$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" property value
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-22T09:24:45+00:00"

# ---------------------------------------------------------------------
# 2. "Mapping Table 1" sheet: correction of the mapping rows for
#    FRAdvanceDirectiveDocument.provision.type and the sourceAttachment
#    sub-elements (id / data / url)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mapping Table 1")

# Row 8: fix source element name (bug fix: .value -> .valueBoolean)
$ws.Range("A8").Value = "FRCDADirectiveAnticipee.valueBoolean"

# Row 12 used to map ".value" -> "sourceAttachment.data"; it now becomes
# the ".id" -> "sourceAttachment.id" mapping
$ws.Range("A12").Value = "FRCDADirectiveAnticipee.entryRelationship.observationMedia.id"
$ws.Range("D12").Value = "FRAdvanceDirectiveDocument.sourceAttachment.id"

# Row 13 keeps its source (...observationMedia.value) but now targets
# "sourceAttachment.data" instead of "sourceAttachment.url"
$ws.Range("D13").Value = "FRAdvanceDirectiveDocument.sourceAttachment.data"

# New row 14: re-add the ".value" -> "sourceAttachment.url" mapping that
# got displaced from row 13. Clone row 13's formatting first, then fill
# in the new values for this row.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A14").Value = "FRCDADirectiveAnticipee.entryRelationship.observationMedia.value"
$ws.Range("C14").Value = "equivalent"
$ws.Range("D14").Value = "FRAdvanceDirectiveDocument.sourceAttachment.url"
